$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data to the right.
$ws.Columns.Item(1).Insert()

# Copy formatting (style) from the old column A (now column B) into the
# newly inserted column A, for the used rows only.
$ws.Range("B1:B10").Copy()
$ws.Range("A1:A10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new "Map viewer name" column with header + friendly names.
# Values are entered in this particular order so that the shared-string
# table is built up in the same sequence as the original authored edit.
$ws.Range("A1").Value = "Map viewer name"
$ws.Range("A2").Value = "Area-based conservation - 30x30 goals"
$ws.Range("A9").Value = "Ecosystem & Natural Process (Re)Creation - opportunities"
$ws.Range("A8").Value = "Ecosystem & Natural Process (Re)Creation - 30x30 goals"
$ws.Range("A10").Value = "Alliance & Partnership Development - land sharing"
$ws.Range("A7").Value = "Site / Area stewardship"
$ws.Range("A3").Value = "Area-based conservation - endemic SAR"
$ws.Range("A4").Value = "Area-based conservation - SAR"
$ws.Range("A5").Value = "Area-based conservation - biodiversity urgency"
$ws.Range("A6").Value = "Area-based conservation - biodiversity opportunities"

# Resize the two leftmost columns to fit the new content (target widths
# 55.7265625 / 66.1796875 characters; the host's ColumnWidth setter only
# resolves to 1/6-character increments, so these are the closest values).
$ws.Columns.Item(1).ColumnWidth = 54.85
$ws.Columns.Item(2).ColumnWidth = 65.35

[void]$ws.Range("A13").Select()
